$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABC-D")

# New rows appended to the "ABC-D" log sheet (recent contest A-Ds).
# Column layout: A=contest(#), B=status, C=date, D=retry?, E=technique, F=comment.
# Note: E87 reuses the sheet's existing "full search" technique entry, whose
# stored text already carries its furigana inline ("全探索" + "ゼンタンサク").
$data = @(
    @(75,  "AC", $null, "二次元累積和", $null),
    @(138, "AC", $null, "dfs", $null),
    @(137, "AC", $null, "貪欲", $null),
    @(136, "AC", $null, $null, $null),
    @(135, "AC", $null, $null, $null),
    @(134, "AC", $true, "全探索ゼンタンサク", "O( N/1 + N/2 + ... + N/N ) = O( N logN )")
)

$startRow = 82
$endRow = $startRow + $data.Count - 1

# Match the existing date-formatted cells (reuse the style, don't mint a new numFmt).
$dateSrc = $ws.Cells.Item($startRow - 1, 3)
$dateRng = $ws.Range("C$startRow`:C$endRow")
$dateSrc.Copy() | Out-Null
$dateRng.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = 43734

    if ($null -ne $vals[2]) {
        $ws.Cells.Item($row, 4).Value = $vals[2]
    }
    if ($null -ne $vals[3]) {
        $ws.Cells.Item($row, 5).Value = $vals[3]
    }
    if ($null -ne $vals[4]) {
        $ws.Cells.Item($row, 6).Value = $vals[4]
    }
}

$ws.Range("F88").Select() | Out-Null
